# Update the "Test Results" column (C) on the LoginData sheet (sheet1) for the
# "Add user" test rows with the latest Selenium failure logs captured after
# re-running the test script (it now gets further, failing at a different
# step than before).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newC2 = @'
no such element: Unable to locate element: {"method":"xpath","selector":"//input[@id='searchUserDetails1']//option[contains(text(),'Distributor')]"}
  (Session info: chrome=136.0.7103.114)
For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Build info: version: '4.30.0', revision: '509c7f17cc*'
System info: os.name: 'Windows 11', os.arch: 'amd64', os.version: '10.0', java.version: '23.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [2ae64c9022fecb543a3b4acb53db6aa1, findElement {value=//input[@id='searchUserDetails1']//option[contains(text(),'Distributor')], using=xpath}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 136.0.7103.114, chrome: {chromedriverVersion: 136.0.7103.113 (76fa3c17824..., userDataDir: C:\Users\NANDAN~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:61320}, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:61320/devtoo..., se:cdpVersion: 136.0.7103.114, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 2ae64c9022fecb543a3b4acb53db6aa1
'@

$newC3 = @'
unexpected alert open: {Alert text : Login failed. Email not found.}
  (Session info: chrome=136.0.7103.114): Login failed. Email not found.
Build info: version: '4.30.0', revision: '509c7f17cc*'
System info: os.name: 'Windows 11', os.arch: 'amd64', os.version: '10.0', java.version: '23.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [5da612ce94c7b2f9efce380f9f77a8fa, getCurrentUrl {}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 136.0.7103.114, chrome: {chromedriverVersion: 136.0.7103.113 (76fa3c17824..., userDataDir: C:\Users\NANDAN~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:61741}, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:61741/devtoo..., se:cdpVersion: 136.0.7103.114, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 5da612ce94c7b2f9efce380f9f77a8fa
'@

$ws.Range("C2").Value = $newC2
$ws.Range("C3").Value = $newC3

$ws.Rows.Item(2).RowHeight = 32.5
$ws.Rows.Item(3).RowHeight = 16.85

Write-Output "Done"
